$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 341, shifting existing
# rows 341-347 down to 343-349 (preserving their data/formatting).
$ws.Rows.Item(341).Resize(2).Insert()

# Populate new row 341
$ws.Cells.Item(341, 1).Value = 5
$ws.Cells.Item(341, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(341, 3).Value = "Maule"
$ws.Cells.Item(341, 4).Value = 44448
$ws.Cells.Item(341, 5).Value = 7
$ws.Cells.Item(341, 6).Value = 100112004
$ws.Cells.Item(341, 7).Value = "Cebolla"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "1a (guarda)"
$ws.Cells.Item(341, 10).Value = 2500
$ws.Cells.Item(341, 11).Value = 2800
$ws.Cells.Item(341, 12).Value = 2800
$ws.Cells.Item(341, 13).Value = 2800
$ws.Cells.Item(341, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(341, 15).Value = "Región del Maule"
$ws.Cells.Item(341, 16).Value = 187
$ws.Cells.Item(341, 17).Value = 15
$ws.Cells.Item(341, 18).Value = "Hortaliza"

# Populate new row 342
$ws.Cells.Item(342, 1).Value = 5
$ws.Cells.Item(342, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(342, 3).Value = "Maule"
$ws.Cells.Item(342, 4).Value = 44448
$ws.Cells.Item(342, 5).Value = 7
$ws.Cells.Item(342, 6).Value = 100112004
$ws.Cells.Item(342, 7).Value = "Cebolla"
$ws.Cells.Item(342, 8).Value = "Sin especificar"
$ws.Cells.Item(342, 9).Value = "1a (guarda)"
$ws.Cells.Item(342, 10).Value = 2500
$ws.Cells.Item(342, 11).Value = 4500
$ws.Cells.Item(342, 12).Value = 4500
$ws.Cells.Item(342, 13).Value = 4500
$ws.Cells.Item(342, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(342, 15).Value = "Región del Maule"
$ws.Cells.Item(342, 16).Value = 180
$ws.Cells.Item(342, 17).Value = 25
$ws.Cells.Item(342, 18).Value = "Hortaliza"
